$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$placeUrl = "https://www.google.it/maps/place/Pantheon/@41.8986108,12.4768729,17z/data=!4m16!1m9!3m8!1s0x132f604f678640a9:0xcad165fa2036ce2c!2sPantheon!8m2!3d41.8986108!4d12.4768729!9m1!1b1!16zL20vMDF4emR6!3m5!1s0x132f604f678640a9:0xcad165fa2036ce2c!8m2!3d41.8986108!4d12.4768729!16zL20vMDF4emR6?entry=ttu"
$giolittiUrl = "https://www.google.it/maps/place/Giolitti/@41.8986032,12.4765859,16.72z/data=!3m1!5s0x132f60519c07e5ef:0x2f3920985ad7eae7!4m16!1m9!3m8!1s0x132f604f678640a9:0xcad165fa2036ce2c!2sPantheon!8m2!3d41.8986108!4d12.4768729!9m1!1b1!16zL20vMDF4emR6!3m5!1s0x132f60519cfa737b:0x38610a40a28f8107!8m2!3d41.9011019!4d12.4771847!16zL20vMGMwMzBu?entry=ttu"

# Replace the URL text in A1 (was the "reviews" URL, now the "place" URL) and
# add a new A2 row with the Giolitti place URL.
$ws.Range("A1").Value = $placeUrl
$ws.Range("A2").Value = $giolittiUrl

# Re-attach the actual hyperlink object on A1 (pointing at the place URL),
# then restore the Hyperlink cell style that Add() always re-stamps.
$ws.Hyperlinks.Add($ws.Range("A1"), $placeUrl, "", "", $placeUrl) | Out-Null
$ws.Range("A1").Style = "Hyperlink"

$ws.Range("E10").Select() | Out-Null
